$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.729.15"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "2.999.86"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.82"
$ws.Range("E5").Value = "  +4.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.56"
$ws.Range("E6").Value = "  +6.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.435"
$ws.Range("E8").Value = "  +4.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.46"
$ws.Range("E9").Value = "  +6.26%  "
$ws.Range("E10").Value = "  +7.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.356"
$ws.Range("E11").Value = "  +2.77%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "3.512.87"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000156"
$ws.Range("E15").Value = "  +12.90%  "
$ws.Range("D16").Value = "56.762.95"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "2.998.42"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.94"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.48"
$ws.Range("E19").Value = "  +4.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.82"
$ws.Range("E20").Value = "  +5.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.85"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.485"
$ws.Range("E23").Value = "  +5.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.22"
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("E25").Value = "  +6.87%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "0.0₃0904"
$ws.Range("E27").Value = "  +8.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.64"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  +7.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +5.67%  "
$ws.Range("E31").Value = "  +7.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.64"
$ws.Range("E32").Value = "  +7.32%  "
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.27"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0679"
$ws.Range("E37").Value = "  +5.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.74"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").Value = "3.034.26"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.05"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "2.281.70"
$ws.Range("E42").Value = "  +7.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.648"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("E44").Value = "  +4.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  +13.02%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.86"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.22"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  +5.97%  "
